# Update cryptocurrency price values (column D) on the active worksheet
# to reflect the latest scraped prices, as produced by the GitHub Actions
# symbol-list updater.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = [ordered]@{
    "D2"  = "274.47"
    "D3"  = "22.94"
    "D4"  = "6.348"
    "D5"  = "0.06236"
    "D7"  = "6.709"
    "D9"  = "0.8323"
    "D11" = "0.1635"
    "D12" = "0.08324"
    "D14" = "0.03116"
    "D15" = "0.09316"
    "D16" = "3.891"
    "D17" = "0.001636"
    "D18" = "0.04770"
    "D19" = "0.006318"
    "D20" = "0.005561"
    "D23" = "3.728"
    "D41" = "0.007027"
    "D42" = "0.1165"
    "D43" = "0.003350"
    "D45" = "0.00006276"
    "D47" = "0.9002"
    "D48" = "0.03386"
}

# The source sheet stores the price column as text (inline strings), so the
# cells must be kept/forced as text rather than letting Excel auto-detect
# these numeric-looking values as Numbers.
foreach ($cell in $updates.Keys) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $updates[$cell]
}
